$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

# Overview sheet: "Latest HO Xliff Generate Date" column G, rows 2-3
$wsOverview.Range("G2").Value = "2017-02-09 10:13:28"
$wsOverview.Range("G3").Value = "2017-02-09 10:13:28"

# zh-cn sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (L)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E3").Value = "mt"
$wsZhCn.Range("H2").Value = "2017-02-09 10:13:07"
$wsZhCn.Range("H3").Value = "2017-02-09 10:13:07"
$wsZhCn.Range("L2").Value = "2017-02-09 10:14:03"
$wsZhCn.Range("L3").Value = "2017-02-09 10:14:03"

# de-de sheet: Priority (E), Correspond Handoff Datetime (H), Correspond Handback DateTime (L)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E3").Value = "mt"
$wsDeDe.Range("H2").Value = "2017-02-09 10:13:28"
$wsDeDe.Range("H3").Value = "2017-02-09 10:13:28"
$wsDeDe.Range("L2").Value = "2017-02-09 10:14:38"
$wsDeDe.Range("L3").Value = "2017-02-09 10:14:38"
